$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.235.40"
$ws.Range("E2").Value = "  -4.96%  "
$ws.Range("D3").Value = "3.252.73"
$ws.Range("E3").Value = "  -7.56%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'592.83"
$ws.Range("E5").Value = "  -4.54%  "
$ws.Range("D6").Value = "'151.05"
$ws.Range("E6").Value = "  -12.73%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.243.83"
$ws.Range("E8").Value = "  -7.76%  "
$ws.Range("D9").Value = "'0.541"
$ws.Range("E9").Value = "  -11.33%  "
$ws.Range("E10").Value = "  -13.57%  "
$ws.Range("D11").Value = "'6.74"
$ws.Range("E11").Value = "  -4.77%  "
$ws.Range("D12").Value = "'0.505"
$ws.Range("E12").Value = "  -13.19%  "
$ws.Range("D13").Value = "'37.99"
$ws.Range("E13").Value = "  -18.12%  "
$ws.Range("D14").Value = "'0.0000242"
$ws.Range("E14").Value = "  -12.62%  "
$ws.Range("D15").Value = "3.773.30"
$ws.Range("E15").Value = "  -7.70%  "
$ws.Range("D16").Value = "67.259.36"
$ws.Range("E16").Value = "  -5.01%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.254.40"
$ws.Range("E17").Value = "  -7.46%  "
$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D18").Value = "'545.05"
$ws.Range("E18").Value = "  -10.57%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").Value = "'0.114"
$ws.Range("E19").Value = "  -6.03%  "
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").Value = "'7.21"
$ws.Range("E20").Value = "  -14.01%  "
$ws.Range("D21").Value = "'15.08"
$ws.Range("E21").Value = "  -14.84%  "
$ws.Range("D22").Value = "'0.761"
$ws.Range("E22").Value = "  -13.80%  "
$ws.Range("D23").Value = "'7.85"
$ws.Range("E23").Value = "  -13.91%  "
$ws.Range("D24").Value = "'85.05"
$ws.Range("E24").Value = "  -13.50%  "
$ws.Range("D25").Value = "'13.51"
$ws.Range("E25").Value = "  -13.45%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").Value = "'3.22"
$ws.Range("E27").Value = "  -13.65%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'29.21"
$ws.Range("E28").Value = "  -13.44%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "'8.01"
$ws.Range("E29").Value = "  -11.53%  "
$ws.Range("D30").Value = "'2.11"
$ws.Range("E30").Value = "  -17.69%  "
$ws.Range("D31").Value = "'2.66"
$ws.Range("E31").Value = "  -11.62%  "
$ws.Range("E32").Value = "  -12.36%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'6.62"
$ws.Range("E33").Value = "  -18.24%  "
$ws.Range("B34").Value = "Bittensor"
$ws.Range("C34").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D34").Value = "'535.49"
$ws.Range("E34").Value = "  -15.42%  "
$ws.Range("D35").Value = "'5.68"
$ws.Range("E35").Value = "  -16.38%  "
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("D37").Value = "'0.0447"
$ws.Range("E37").Value = "  -7.17%  "
$ws.Range("D38").Value = "'53.46"
$ws.Range("E38").Value = "  -5.83%  "
$ws.Range("D39").Value = "'0.0852"
$ws.Range("E39").Value = "  -14.60%  "
$ws.Range("D40").Value = "'9.14"
$ws.Range("E40").Value = "  -15.29%  "
$ws.Range("E41").Value = "  -11.45%  "
$ws.Range("D42").Value = "2.929.07"
$ws.Range("E42").Value = "  -12.51%  "
$ws.Range("E43").Value = "  -22.61%  "
$ws.Range("E44").Value = "  -16.69%  "
$ws.Range("D45").Value = "0.0₃0576"
$ws.Range("E45").Value = "  -19.78%  "
$ws.Range("D46").Value = "'26.38"
$ws.Range("E46").Value = "  -16.93%  "
$ws.Range("E47").Value = "  -15.55%  "
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "'127.08"
$ws.Range("E49").Value = "  -5.28%  "
$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").Value = "'2.34"
$ws.Range("E50").Value = "  -20.53%  "
$ws.Range("D51").Value = "'0.113"
$ws.Range("E51").Value = "  -12.68%  "
